$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 434; this shifts the existing rows 434:520
# down to 436:522 (and Excel's used range / dimension grows to R522
# automatically).
$ws.Rows("434:435").Insert()

# The two inserted rows share every "boilerplate" column with the rest of
# this Mercado/Categoria block, so copy that straight from the row that
# used to be 434 (now pushed down to row 436).
$ws.Range("A434:C435").Value2 = $ws.Range("A436:C437").Value2
$ws.Range("E434:H435").Value2 = $ws.Range("E436:H437").Value2
$ws.Range("N434:O435").Value2 = $ws.Range("N436:O437").Value2
$ws.Range("Q434:R435").Value2 = $ws.Range("Q436:R437").Value2

# New "Primera" quality row for the new (latest) date.
$ws.Range("D434").Value2 = 45211
$ws.Range("I434").Value2 = "Primera"
$ws.Range("J434").Value2 = 1600
$ws.Range("K434").Value2 = 500
$ws.Range("L434").Value2 = 600
$ws.Range("M434").Value2 = 550
$ws.Range("P434").Value2 = 183

# New "Segunda" quality row for the new (latest) date.
$ws.Range("D435").Value2 = 45211
$ws.Range("I435").Value2 = "Segunda"
$ws.Range("J435").Value2 = 1300
$ws.Range("K435").Value2 = 400
$ws.Range("L435").Value2 = 450
$ws.Range("M435").Value2 = 425
$ws.Range("P435").Value2 = 142
